$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# Cells whose new value looks like a plain number get an explicit Text
# number format first, so Excel stores the exact string instead of
# silently converting it to a floating-point number (e.g. "141.10" -> 141.1).

$ws.Range("D2").Value = "57.780.38"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "3.109.36"
$ws.Range("E3").Value = "  +1.09%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "524.63"
$ws.Range("E5").Value = "  +1.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.10"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "3.107.49"
$ws.Range("E8").Value = "  +1.13%  "
$ws.Range("E9").Value = "  +0.48%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.28"
$ws.Range("E10").Value = "  -0.13%  "
$ws.Range("E11").Value = "  +0.80%  "
$ws.Range("E12").Value = "  +2.28%  "
$ws.Range("D13").Value = "3.645.93"
$ws.Range("E13").Value = "  +1.18%  "
$ws.Range("E14").Value = "  +1.77%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.23"
$ws.Range("E15").Value = "  +2.09%  "
$ws.Range("E16").Value = "  +0.69%  "
$ws.Range("D17").Value = "57.888.78"
$ws.Range("E17").Value = "  +0.38%  "
$ws.Range("D18").Value = "3.108.83"
$ws.Range("E18").Value = "  +0.99%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.11"
$ws.Range("E19").Value = "  +0.72%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.82"
$ws.Range("E20").Value = "  -1.54%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.06"
$ws.Range("E21").Value = "  -0.60%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "337.34"
$ws.Range("E22").Value = "  +0.85%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.507"
$ws.Range("E24").Value = "  +1.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "66.59"
$ws.Range("E25").Value = "  +0.95%  "
$ws.Range("E26").Value = "  -1.00%  "
$ws.Range("E27").Value = "  +0.19%  "
$ws.Range("D28").Value = "0.0₃0928"
$ws.Range("E28").Value = "  +1.68%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.57"
$ws.Range("E29").Value = "  +3.35%  "
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.23"
$ws.Range("E31").Value = "  +0.44%  "
$ws.Range("E32").Value = "  +2.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.20"
$ws.Range("E33").Value = "  +2.56%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "20.94"
$ws.Range("E34").Value = "  +0.22%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "154.54"
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.62"
$ws.Range("E36").Value = "  +3.41%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.07"
$ws.Range("E37").Value = "  +2.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "27.16"
$ws.Range("E38").Value = "  -0.04%  "
$ws.Range("E39").Value = "  +1.50%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0667"
$ws.Range("E40").Value = "  -1.42%  "
$ws.Range("D41").Value = "3.154.67"
$ws.Range("E41").Value = "  +1.26%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.683"
$ws.Range("E42").Value = "  +3.67%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.90"
$ws.Range("E43").Value = "  -0.70%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "36.94"
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.999"
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.46"
$ws.Range("E46").Value = "  +5.81%  "
$ws.Range("D47").Value = "2.281.02"
$ws.Range("E47").Value = "  +0.67%  "
$ws.Range("E48").Value = "  +0.57%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.984"
$ws.Range("E49").Value = "  +5.85%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.66"
$ws.Range("E50").Value = "  +1.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.00"
$ws.Range("E51").Value = "  +2.08%  "
